# Update "paises" (countries) workbook: refresh case numbers and the
# "last updated" timestamp, matching a newer data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" banner text (row 1) -----------------------
$ws.Range("A1").Value = "Datos actualizados a 22 de Abril de 2020 a las 18:22"

# --- Row 6 (Italia) --------------------------------------------------------
$ws.Range("B6").Value = 187327
$ws.Range("C6").Value = 3370
$ws.Range("D6").Value = 54543
$ws.Range("E6").Value = 107699
$ws.Range("F6").Value = 2384
$ws.Range("G6").Value = 437
$ws.Range("H6").Value = 25085

# --- Row 16 (Canada) --------------------------------------------------------
$ws.Range("B16").Value = 38967
$ws.Range("C16").Value = 545
$ws.Range("D16").Value = 13647
$ws.Range("E16").Value = 23447
$ws.Range("G16").Value = 39
$ws.Range("H16").Value = 1873

# --- Row 18 (Suiza) --------------------------------------------------------
$ws.Range("E18").Value = 7359
$ws.Range("G18").Value = 31
$ws.Range("H18").Value = 1509

# --- Row 40 (Noruega) -------------------------------------------------------
$ws.Range("E40").Value = 7057
$ws.Range("G40").Value = 4
$ws.Range("H40").Value = 186

# --- Row 58/59: Argelia overtakes Tailandia in rank -------------------------
# Argelia now has more total cases than Tailandia, so Argelia moves to
# row 58 (rank 58) and Tailandia drops to row 59 (rank 59).
$ws.Range("A58").Value = "Argelia"
$ws.Range("B58").Value = 2910
$ws.Range("C58").Value = 99
$ws.Range("D58").Value = 1204
$ws.Range("E58").Value = 1304
$ws.Range("F58").Value = 40
$ws.Range("G58").Value = 10
$ws.Range("H58").Value = 402

$ws.Range("A59").Value = "Tailandia"
$ws.Range("B59").Value = 2826
$ws.Range("C59").Value = 15
$ws.Range("D59").Value = 2352
$ws.Range("E59").Value = 425
$ws.Range("F59").Value = 61
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = 49

# --- Row 61 (Grecia) --------------------------------------------------------
$ws.Range("F61").Value = 55

# --- Row 118 (Isla de Man) ---------------------------------------------------
$ws.Range("D118").Value = 212
$ws.Range("E118").Value = 80
$ws.Range("F118").Value = 20
$ws.Range("G118").Value = 6
$ws.Range("H118").Value = 15

# --- Row 153 (Cabo Verde) ----------------------------------------------------
$ws.Range("B153").Value = 73
$ws.Range("C153").Value = 5
$ws.Range("E153").Value = 71

# --- Row 189 (San Cristobal y Nieves) ---------------------------------------
$ws.Range("D189").Value = 1
$ws.Range("E189").Value = 14
